# Applies the "Reports" sheet addition + related view-state changes
# described by the commit diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add a new worksheet named "Reports" after the last existing sheet
#    (so it becomes sheet #3, after "Sheet1" and "MIGRACIONES DE
#    FRAMEWORK").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$reports = $wb.Worksheets.Add($null, $lastSheet)
$reports.Name = "Reports"

# ---------------------------------------------------------------------
# 2. Fill in the new sheet's content (matches the new shared-string
#    entries / cell map from the diff).
# ---------------------------------------------------------------------
$reports.Range("C2").Value = "Para trabajar con RDLC y Reportes es necesario tomar en cuenta lo siguiente"

$reports.Range("C3").Value = "1)"
$reports.Range("D3").Value = "Hasta feb del 2024 ReportViewer no funcionaba con .NetCore por lo que fue necesario primero crear un nuevo proyecto"

$reports.Range("D4").Value = "Que tenga un ambiente para diseñar el reporte y esto lo conseguimos con WindowsForm App (.NET Framework) con .Net Framework 4.7.2"

$reports.Range("D5").Value = "Como referencia el video "
$reports.Range("G5").Value = "https://www.youtube.com/watch?v=41RcaFPphTA"

$reports.Range("C6").Value = "2) "
$reports.Range("D6").Value = "Bajar Microsoft RDLC Report Designer si es que no lo tenemos"

# ---------------------------------------------------------------------
# 3. Match the recorded selection on the new sheet (D7).
# ---------------------------------------------------------------------
[void]$reports.Range("D7").Select()

# ---------------------------------------------------------------------
# 4. Scroll the "MIGRACIONES DE FRAMEWORK" sheet's view down (its
#    recorded topLeftCell moves to A16) while leaving its selection at
#    C4. It is no longer the active/tabSelected sheet after this edit.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MIGRACIONES DE FRAMEWORK")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# 5. Re-activate the new "Reports" sheet so it ends up as the
#    workbook's active tab (activeTab goes from 1 -> 2) and carries
#    tabSelected="1" on save.
# ---------------------------------------------------------------------
$reports.Activate()
